# Adds "wrong answer" option columns (Valik1/Valik2/Valik3) to the two
# question sheets, pushing the existing "Vanusegrupp" column out to F.
# Cell writes are ordered to match the original authoring session so the
# shared-string table grows in the same sequence.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Grupp_1")
$ws2 = $wb.Worksheets.Item("Grupp_2")

# --- Grupp_1 (sheet1): insert 3 columns before the old column C --------
$ws1.Range("C1:E1").EntireColumn.Insert()

$ws1.Range("C1").Value = "Valik1"
$ws1.Range("D1").Value = "Valik2"
$ws1.Range("E1").Value = "Vlaik3"

$ws1.Range("C2").Value = 42
$ws1.Range("D2").Value = "uni"
$ws1.Range("E2").Value = "õlu"

$ws1.Range("C3").Value = 2
$ws1.Range("D3").Value = 3
$ws1.Range("E3").Value = 4

# --- Grupp_2 (sheet2): insert 3 columns before the old column C --------
$ws2.Range("C1:E1").EntireColumn.Insert()

$ws2.Columns("C").ColumnWidth = 23.85546875
$ws2.Columns("D:E").ColumnWidth = 20.28515625

$ws2.Range("C1").Value = "Valik1"
$ws2.Range("D1").Value = "Valik2"
$ws2.Range("E1").Value = "valik3"

$ws2.Range("C2").Value = "Luts"
$ws2.Range("D2").Value = "Köstrihärra"
$ws2.Range("E2").Value = "Toots"

$ws2.Range("C3").Value = "Silinder"
$ws2.Range("D3").Value = "Trapets"
$ws2.Range("E3").Value = "Võrdhaarne trapets"

$ws2.Range("C4").Value = "Gravitatsioonikiirendus"
$ws2.Range("D4").Value = "Kliirensikiirendus"
$ws2.Range("F4").Value = 2

$ws2.Range("C5").Value = "Keskringjoone"
$ws2.Range("D5").Value = "Puutuja "
$ws2.Range("E5").Value = "Keskristsirge"

$ws2.Range("C6").Value = "Maksim"
$ws2.Range("D6").Value = "Einstein"
$ws2.Range("E6").Value = "Galileo"
$ws2.Range("F6").Value = 2

# E4 ("Kiirendusvõistlus") was typed last of the new distractor strings.
$ws2.Range("E4").Value = "Kiirendusvõistlus"

# --- Selection / active sheet changes -----------------------------------
[void]$ws1.Range("C4").Select()
[void]$ws2.Activate()
[void]$ws2.Range("A17").Select()
